$d = $word.ActiveDocument

# -----------------------------------------------------------------
# Change 1: remove the _GoBack bookmark that wraps "pyxlinks " (it
# sat between the "and " run and the "directories to Chimera share
# directory. The " run). Removing the bookmark leaves two adjacent
# runs with identical formatting (<w:rFonts w:cs="Times New Roman"/>)
# - " " and "directories to Chimera share directory. The " - which
# Word coalesces into a single run once either of them is re-written.
# -----------------------------------------------------------------
$bm = $d.Bookmarks.Item("_GoBack")
$bmRange = $bm.Range
$bmEnd = $bmRange.End
$bm.Delete()

$mergeFind = $d.Content
$mergeFind.Find.Execute("directories to Chimera share directory. The ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$mergeEnd = $mergeFind.End

$mergeRange = $d.Range($bmEnd - 1, $mergeEnd)
$mergeOriginal = $mergeRange.Text
# touch the text (append then restore) to force the run-coalescing
# pass to fold the two runs into one, without altering the content
$mergeRange.Text = $mergeOriginal + "~"
$mergeRange2 = $d.Range($bmEnd - 1, $mergeEnd + 1)
$mergeRange2.Text = $mergeOriginal

# -----------------------------------------------------------------
# Change 2: "Setup -> Load project" becomes two italic runs:
# "Menu File" and " -> Load project".
# -----------------------------------------------------------------
$f2 = $d.Content
$f2.Find.Execute("Setup -> Load project", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$start2 = $f2.Start
$end2 = $f2.End

$firstPart = $d.Range($start2, $start2 + 5)   # "Setup"
$firstPart.Text = "Menu File"

$secondStart = $start2 + 9                    # length of "Menu File"
$secondLen = $end2 - ($start2 + 5)             # length of " -> Load project"
$secondPart = $d.Range($secondStart, $secondStart + $secondLen)
# toggle a formatting property to keep this as its own run instead of
# re-merging with the previous (now textually-identical-style) run
$secondPart.Bold = 1
$secondPart.Bold = 0

# -----------------------------------------------------------------
# Change 3: "Subunits tab -> select 4C3H.pdb -> Color all subunits"
# becomes "Subunits tab -> select 4C3H.pdb -> Color all " (italic)
# followed by a new bookmark "_GoBack" wrapping a plain run "button".
# -----------------------------------------------------------------
$f3 = $d.Content
$f3.Find.Execute("Subunits tab -> select 4C3H.pdb -> Color all subunits", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$start3 = $f3.Start
$end3 = $f3.End

$tail = $d.Range($end3 - 8, $end3)             # "subunits"
$tail.Text = ""

$insertPoint = $end3 - 8
$buttonRange = $d.Range($insertPoint, $insertPoint)
$buttonRange.InsertAfter("button")

$newBmRange = $d.Range($insertPoint, $insertPoint + 6)   # "button"
$d.Bookmarks.Add("_GoBack", $newBmRange)
